$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2024-04-28 Sunday" "2024-04-29 Monday"

Replace-Text "979×6=5874" "822×2=1644"
Replace-Text "884×4=3536" "669×7=4683"
Replace-Text "773×3=2319" "850×2=1700"
Replace-Text "814×2=1628" "983×8=7864"
Replace-Text "818×5=4090" "143×4=572"
Replace-Text "939×3=2817" "139×4=556"
Replace-Text "769×2=1538" "213×5=1065"
Replace-Text "419×9=3771" "627×8=5016"
Replace-Text "471×8=3768" "307×2=614"
Replace-Text "640×5=3200" "586×7=4102"
Replace-Text "757×7=5299" "518×6=3108"
Replace-Text "709×6=4254" "371×6=2226"
Replace-Text "293×4=1172" "502×3=1506"
Replace-Text "175×3=525" "268×8=2144"
Replace-Text "236×2=472" "506×3=1518"
Replace-Text "850×7=5950" "245×4=980"
Replace-Text "766×3=2298" "954×2=1908"
Replace-Text "297×5=1485" "988×9=8892"
Replace-Text "671×6=4026" "782×5=3910"
Replace-Text "732×2=1464" "267×3=801"
Replace-Text "405×2=810" "823×6=4938"
Replace-Text "433×4=1732" "836×7=5852"
Replace-Text "489×6=2934" "394×6=2364"
Replace-Text "164×6=984" "825×9=7425"
Replace-Text "538×4=2152" "570×6=3420"
